$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill_Info_List")
$ws.Activate()

# Rows where column K changes from 0 to 1 (cool time 1 for hit effect)
$rows = @(4..51) + @(54..117) + @(119..136)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 11).Value = 1
}

# Update the view: select K118, then scroll so row 109 is at the top of the window
$ws.Range("K118").Select()
$appWin = $excel.ActiveWindow
$appWin.ScrollRow = 109
$appWin.ScrollColumn = 1
